# Actualización automática 2025-09-24 15:30:09
#
# This workbook is a static report (no live formulas) with three sheets:
#   "VENTAS POR GRUPO"      - sales by product group per client
#   "VENTA MENSUAL"         - sales by month per client
#   "CUMPLIMIENTO MENSUAL"  - budget vs. sales compliance per product group
#
# A new "septiembre" sale was recorded for advisor ALMEIDA CUATIN JHONATHANN
# CARLOS: client COMFALASDI (PORCELANATO, 2728.76) and client SANCHEZ
# SARMIENTO ANDRES FERNANDO (FREGADEROS DE COCINA 143.42 + INODOROS 805.5).
# This ripples into the monthly totals, the "x de 33" non-zero counters, and
# the monthly compliance (PRESUPUESTO/VENTA/POR CUMPLIR/CUMPLIMIENTO) rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": per-client sale amounts by product group
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M8").Value = 2728.76    # COMFALASDI ... -> PORCELANATO
$wsGrupo.Range("E27").Value = 143.42    # SANCHEZ SARMIENTO ... -> FREGADEROS DE COCINA
$wsGrupo.Range("H27").Value = 805.5     # SANCHEZ SARMIENTO ... -> INODOROS

# Row 35 footer counts how many of the 33 data rows are non-zero per column
$wsGrupo.Range("E35").Value = "2 de 33"
$wsGrupo.Range("H35").Value = "1 de 33"
$wsGrupo.Range("M35").Value = "8 de 33"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": per-client sale amounts by month
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F8").Value = 2728.76   # COMFALASDI ... -> septiembre
$wsMensual.Range("F27").Value = 948.92   # SANCHEZ SARMIENTO ... -> septiembre
$wsMensual.Range("F35").Value = 21441.41 # septiembre total

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": budget vs. sales per product group
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 4: FREGADEROS DE COCINA
$wsCumpl.Range("D4").Value = 289.37
$wsCumpl.Range("E4").Value = 232.24144263264
$wsCumpl.Range("F4").Value = 0.5547616028887565

# Row 6: INODOROS
$wsCumpl.Range("D6").Value = 805.5
$wsCumpl.Range("E6").Value = 8.623430808872968
$wsCumpl.Range("F6").Value = 0.9894077108181186

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 13940.21
$wsCumpl.Range("E12").Value = 8493.5453751766
$wsCumpl.Range("F12").Value = 0.6213944017337873

# Row 15: TOTAL
$wsCumpl.Range("D15").Value = 21698.64
$wsCumpl.Range("E15").Value = 17044.37881339592
$wsCumpl.Range("F15").Value = 0.5600658044875274
